# Generate Report for Handback
# Updates the localization-status report:
#  - zh-cn / de-de sheets: Status -> "Handed back: in sync with en-US",
#    refresh the "Latest Handback DateTime" stamp, and clear the stale
#    "Error Detail" warning now that the handback is in sync.
#  - Widen the "Status" columns (and shrink the "Error Detail" column) so
#    the new, longer status text / now-empty error column read cleanly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet (status roll-up columns for each language) --------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet -----------------------------------------------------
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-31 15:04:09"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet -------------------------------------------------------
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-31 15:04:33"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments -----------------------------------------
# Overview: widen the zh-cn / de-de status columns (E, F)
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: widen the Status column (C) and shrink Error Detail (P)
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
